$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete L:N (contrato/modalidade/registro) and O:Q columns content
# that no longer exist in the updated layout.
$ws.Range("O1:Q3").Clear()

# --- Row 2 ---
$ws.Range("A2").Value = "417823 - PREMIUM SAÚDE S.A"
$ws.Range("B2").Value = "19-04-2023"
$ws.Range("C2").Value = 45026.38854166667
$ws.Range("C2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D2").Value = 12161100
$ws.Range("E2").Value = 8597675
$ws.Range("F2").Value = "CAIO HENRIQUE RODRIGUES FERNANDES"
$ws.Range("G2").Value = 17042707664
$ws.Range("H2").Value = "Interlocutora que se identifica como mãe do beneficiário, questiona a falta de atendimento para Consulta com Otorrinolaringologista, Audiometria Tonal e Vocal, Impedanciometria, Videoendoscopia nasossinusal. A solicitação foi feita à Operadora no dia 23/03/2023, para realização no município BETIM. A operadora não apresenta resposta ao pedido, informou apenas que buscaria profissional. Protocolo: 3682532023041042947 - Data: 10/04/2023."
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = "NO"
$ws.Range("K2").Value = "Assistencial"
$ws.Range("L2").Value = "Cláudio"
$ws.Range("M2").Value = "Vieira"
$ws.Range("N2").Value = "Amantino"

# --- Row 3 ---
$ws.Range("A3").Value = "417823 - PREMIUM SAÚDE S.A"
$ws.Range("B3").Value = "19-04-2023"
$ws.Range("C3").Value = 45026.44988425926
$ws.Range("C3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D3").Value = 12161386
$ws.Range("E3").Value = 8598016
$ws.Range("F3").Value = "DIEGO SANTOS DE ALMEIDA"
$ws.Range("G3").Value = 1990602665
$ws.Range("H3").Value = "Interlocutora, que se identifica como esposa  do beneficiário, questiona a falta de atendimento para cartão com orçamento e valores dos matérias para cirurgia  reconstrução do ligamento cruzado anterior  . A solicitação foi feita à Operadora em julho/2022, para realização no município patos de minas . A operadora  não apresenta resposta ao pedido, operadora já lhe encaminhou uma carta mas não esta\ com os valores dos matérias com o hospital deseja    Protocolo:36825320230410426644 data:1/04/2023."
$ws.Range("I3").Value = 4
$ws.Range("J3").Value = "NO"
$ws.Range("K3").Value = "Assistencial"
$ws.Range("L3").Value = "João"
$ws.Range("M3").Value = "Silva"
$ws.Range("N3").Value = "Batista"
